# Applies the "Export with no is_pref and no lev distance" re-export:
# updates columns B (id) and C (speaker_variant) with the new row
# ordering/values, and clears column D (is_prefered) for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "#tanne"
$ws.Cells.Item(2, 3).Value = "Tanne"
$ws.Cells.Item(2, 4).Value = ""

$ws.Cells.Item(3, 2).Value = "#neel"
$ws.Cells.Item(3, 3).Value = "Neel"
$ws.Cells.Item(3, 4).Value = ""

$ws.Cells.Item(4, 2).Value = "#rog"
$ws.Cells.Item(4, 3).Value = "Rog"
$ws.Cells.Item(4, 4).Value = ""

$ws.Cells.Item(5, 2).Value = "#ian-sal"
$ws.Cells.Item(5, 3).Value = "Ian Sal"
$ws.Cells.Item(5, 4).Value = ""

$ws.Cells.Item(6, 2).Value = "#ian.-mart"
$ws.Cells.Item(6, 3).Value = "Ian. Mart"
$ws.Cells.Item(6, 4).Value = ""

$ws.Cells.Item(7, 2).Value = "#ioosjen"
$ws.Cells.Item(7, 3).Value = "Ioosjen"
$ws.Cells.Item(7, 4).Value = ""

$ws.Cells.Item(8, 2).Value = "#trijn"
$ws.Cells.Item(8, 3).Value = "Trijn"
$ws.Cells.Item(8, 4).Value = ""

$ws.Cells.Item(9, 2).Value = "#griet"
$ws.Cells.Item(9, 3).Value = "Griet"
$ws.Cells.Item(9, 4).Value = ""

$ws.Cells.Item(10, 2).Value = "#kennis"
$ws.Cells.Item(10, 3).Value = "Kennis"
$ws.Cells.Item(10, 4).Value = ""

$ws.Cells.Item(11, 2).Value = "#mart"
$ws.Cells.Item(11, 3).Value = "Mart"
$ws.Cells.Item(11, 4).Value = ""

$ws.Cells.Item(12, 2).Value = "#elsje"
$ws.Cells.Item(12, 3).Value = "Elsje"
$ws.Cells.Item(12, 4).Value = ""

$ws.Cells.Item(13, 2).Value = "#cornel"
$ws.Cells.Item(13, 3).Value = "Cornel"
$ws.Cells.Item(13, 4).Value = ""

$ws.Cells.Item(14, 2).Value = "#ioost"
$ws.Cells.Item(14, 3).Value = "Ioost"
$ws.Cells.Item(14, 4).Value = ""

$ws.Cells.Item(15, 2).Value = "#lubb"
$ws.Cells.Item(15, 3).Value = "Lubb"
$ws.Cells.Item(15, 4).Value = ""

$ws.Cells.Item(16, 2).Value = "#olyvier"
$ws.Cells.Item(16, 3).Value = "Olyvier"
$ws.Cells.Item(16, 4).Value = ""

$ws.Cells.Item(17, 2).Value = "#rogier"
$ws.Cells.Item(17, 3).Value = "Rogier"
$ws.Cells.Item(17, 4).Value = ""

$ws.Cells.Item(18, 2).Value = "#jan-sal"
$ws.Cells.Item(18, 3).Value = "Jan Sal"
$ws.Cells.Item(18, 4).Value = ""

$ws.Cells.Item(19, 2).Value = "#lubbert"
$ws.Cells.Item(19, 3).Value = "Lubbert"
$ws.Cells.Item(19, 4).Value = ""

$ws.Cells.Item(20, 2).Value = "#ioos"
$ws.Cells.Item(20, 3).Value = "Ioos"
$ws.Cells.Item(20, 4).Value = ""

$ws.Cells.Item(21, 2).Value = "#marcel"
$ws.Cells.Item(21, 3).Value = "Marcel"
$ws.Cells.Item(21, 4).Value = ""

$ws.Cells.Item(22, 2).Value = "#trijn-ra"
$ws.Cells.Item(22, 3).Value = "Trijn Ra"
$ws.Cells.Item(22, 4).Value = ""

$ws.Cells.Item(23, 2).Value = "#tan"
$ws.Cells.Item(23, 3).Value = "Tan"
$ws.Cells.Item(23, 4).Value = ""

$ws.Cells.Item(24, 2).Value = "#ian"
$ws.Cells.Item(24, 3).Value = "Ian"
$ws.Cells.Item(24, 4).Value = ""

$ws.Cells.Item(25, 2).Value = "#elsjen"
$ws.Cells.Item(25, 3).Value = "Elsjen"
$ws.Cells.Item(25, 4).Value = ""

$ws.Cells.Item(26, 2).Value = "#iuriaen"
$ws.Cells.Item(26, 3).Value = "Iuriaen"
$ws.Cells.Item(26, 4).Value = ""

